$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the existing sheets ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Copy header style (bold, bordered, centered) from an existing header cell
# so the new header cells reuse the same cell style rather than creating
# duplicate style entries.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Copy the date-formatted style from an existing data cell in column A
# down the new column A so the dates render the same way.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A75").PasteSpecial(-4122)

$rows = @(
    "44990.99999999999,176,-221.2937810430981,574.3359843326757"
    "45011.99999999999,179,-226.3821839812109,595.6056612962052"
    "45018.99999999999,180,-222.3664907046791,546.7712207290098"
    "45032.99999999999,182,-205.370848507278,583.9405748989474"
    "45039.99999999999,184,-224.5140567050223,586.9463346520793"
    "45046.99999999999,185,-183.3201383738543,590.6259359834081"
    "45053.99999999999,186,-212.62761703661,626.9786853475017"
    "45067.99999999999,188,-199.4493316050259,576.9090993118384"
    "45074.99999999999,189,-204.1951766372016,598.9840531078722"
    "45081.99999999999,190,-182.5575309585291,596.6444608035571"
    "45095.99999999999,192,-199.1057167160041,626.5821228067858"
    "45102.99999999999,194,-243.5021326803785,602.6143746940692"
    "45109.99999999999,195,-203.9038894116246,592.3105396825483"
    "45116.99999999999,196,-220.2742139471194,594.3140026458475"
    "45123.99999999999,197,-200.4364900380933,596.557983543589"
    "45130.99999999999,198,-173.3574446461157,569.8499079520102"
    "45137.99999999999,199,-189.2585270603783,579.7096157606869"
    "45151.99999999999,201,-229.9664639681319,600.0313160607567"
    "45158.99999999999,202,-178.2481797757084,609.1954842476508"
    "45165.99999999999,203,-185.9090853129998,624.1073957794678"
    "45172.99999999999,205,-201.7310784401631,613.7510651393967"
    "45179.99999999999,206,-214.1130604171393,618.7959350098621"
    "45193.99999999999,208,-179.2694108598149,601.1561200387611"
    "45200.99999999999,209,-181.1485187884218,566.6632546365489"
    "45207.99999999999,210,-163.3569947326512,584.4586643572399"
    "45221.99999999999,212,-177.2956392054644,619.5052351942062"
    "45228.99999999999,213,-184.4217862618586,613.2419314558463"
    "45235.99999999999,214,-194.5374279525589,604.126465115535"
    "45242.99999999999,216,-187.2839303501831,649.4015600678898"
    "45249.99999999999,217,-211.1981585782617,600.937872278896"
    "45256.99999999999,218,-170.5570822811372,622.6510501428364"
    "45263.99999999999,219,-173.6956419490753,622.2609869189719"
    "45270.99999999999,220,-184.1796549656065,623.6169693120794"
    "45277.99999999999,221,-165.0314081308447,637.3423953766567"
    "45298.99999999999,224,-160.973160292597,645.5603295343711"
    "45305.99999999999,226,-183.3617361971575,625.4814051921632"
    "45312.99999999999,227,-155.8633948847795,593.471461839224"
    "45340.99999999999,231,-155.636621029931,631.6883088334504"
    "45347.99999999999,232,-168.5643246773315,641.3419824957049"
    "45354.99999999999,233,-156.2883394883678,623.5455756256408"
    "45361.99999999999,234,-165.2699240195543,612.0662539262023"
    "45382.99999999999,238,-153.7323281127556,666.486186851071"
    "45389.99999999999,239,-182.6668810794226,630.9336347633169"
    "45396.99999999999,240,-147.6301595599797,629.6331569965654"
    "45403.99999999999,241,-177.1872050064341,647.2451344780164"
    "45410.99999999999,242,-186.9074291616142,628.9249059077263"
    "45417.99999999999,243,-159.2862493158628,665.0881067124272"
    "45424.99999999999,244,-115.1312689617555,639.6232080673028"
    "45431.99999999999,245,-117.6462180060578,638.1051525229697"
    "45438.99999999999,247,-145.408855490803,646.2959143681497"
    "45459.99999999999,250,-165.3965932120313,638.6484102676289"
    "45466.99999999999,251,-121.3984117594224,684.6484196573487"
    "45487.99999999999,254,-123.0730396216527,636.2621775965763"
    "45501.99999999999,256,-141.4407934007868,674.8099136923742"
    "45508.99999999999,258,-132.369898481408,649.5234965979023"
    "45515.99999999999,259,-152.9410819079981,662.4296389337251"
    "45529.99999999999,261,-127.0915982294549,665.5596720784952"
    "45536.99999999999,262,-141.4127289914873,654.7433161650431"
    "45543.99999999999,263,-136.1247447537564,669.8993195579857"
    "45557.99999999999,265,-142.173890175107,635.7035723261953"
    "45564.99999999999,266,-93.42307114959921,680.931514958421"
    "45571.99999999999,267,-138.9702489399738,656.658163598227"
    "45578.99999999999,269,-127.3100267879912,669.6100712240603"
    "45585.99999999999,270,-134.3584169603437,645.854092674898"
    "45599.99999999999,272,-150.9920962297153,670.242770598588"
    "45634.99999999999,277,-115.610916849816,676.0798809664109"
    "45641.99999999999,279,-111.3943753228034,647.2118612470805"
    "45648.99999999999,280,-121.6326832750753,701.1548263910462"
    "45655.99999999999,281,-112.0773823057843,673.2747315810924"
    "45662.99999999999,282,-175.2527464688468,643.3669445737552"
    "45669.99999999999,283,-131.3319248572127,664.6098885196947"
    "45676.99999999999,284,-144.3483156169267,666.6447850652382"
    "45683.99999999999,285,-132.7123081849538,672.1591387560208"
    "45690.99999999999,286,-132.5266567503968,654.5950616108041"
)

$r = 2
foreach ($line in $rows) {
    $parts = $line.Split(",")
    $wsForecast.Cells.Item($r, 1).Value = [double]$parts[0]
    $wsForecast.Cells.Item($r, 2).Value = [double]$parts[1]
    $wsForecast.Cells.Item($r, 3).Value = [double]$parts[2]
    $wsForecast.Cells.Item($r, 4).Value = [double]$parts[3]
    $r++
}
